$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1374
$ws1.Range("F4").Value = 88

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1374
$ws4.Range("F4").Value = 88
